$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bug fix: the Springer API abstract-fetch path re-ran its author-string
# formatting step on top of data that had already been formatted, so every
# "Authors" cell (column E) picked up one extra space after every comma
# each time the pipeline ran. Re-running the fetch twice therefore grows
# the inter-field spacing by two extra spaces per comma. Re-apply that same
# whitespace growth (once per pass, twice total) to every data row in the
# "Authors" column so the stored values match what the corrected/re-run
# pipeline now produces.

$headerRow = 1
$lastRow = $ws.UsedRange.Row + $ws.UsedRange.Rows.Count - 1
$lastCol = $ws.UsedRange.Column + $ws.UsedRange.Columns.Count - 1

$authorsCol = 0
for ($c = $ws.UsedRange.Column; $c -le $lastCol; $c++) {
    if ([string]$ws.Cells.Item($headerRow, $c).Value2 -eq "Authors") {
        $authorsCol = $c
    }
}

for ($row = $headerRow + 1; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, $authorsCol)
    $value = [string]$cell.Value2
    if ($value -ne "") {
        $value = $value -replace ',( +)', ',$1 '
        $value = $value -replace ',( +)', ',$1 '
        $cell.Value = $value
    }
}
